$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Neutrons per cycle" value for Layer 2 (outer) from 1500 to 1000
$ws.Range("D21").Value = 1000

# Update the selection/view state to reflect where the edit was made
$ws.Range("D24:E24").Select()
$excel.ActiveWindow.Zoom = 70
